$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Advance task Nunit Project" data entry — append new rows of test
# data to the Language, Skill, Education and Certificate lookup sheets.
# ---------------------------------------------------------------------

$wsLanguage    = $wb.Worksheets.Item("Language")
$wsSkill       = $wb.Worksheets.Item("Skill")
$wsEducation   = $wb.Worksheets.Item("Education")
$wsCertificate = $wb.Worksheets.Item("Certificate")

# --- Language sheet: new rows 3 & 4 (row 5 filled in further below,
#     after Skill/Certificate, to match the original authoring order) --
$wsLanguage.Range("A3").Value = "Gujarati"
$wsLanguage.Range("B3").Value = "Hindi"
$wsLanguage.Range("A4").Value = "Marathi"
$wsLanguage.Range("B4").Value = "Gujarati"

# --- Skill sheet: new rows 3 & 4 ---------------------------------------
$wsSkill.Range("A3").Value = "Reading"
$wsSkill.Range("B3").Value = "Painting"
$wsSkill.Range("A4").Value = "Painting"
$wsSkill.Range("B4").Value = "Speaking"

# --- Certificate sheet: new row 3 --------------------------------------
$wsCertificate.Range("A3").Value = "Google"
$wsCertificate.Range("B3").Value = "Australia"
$wsCertificate.Range("C3").Value = "ISTQB2"
$wsCertificate.Range("D3").Value = "India"

# --- back to Language sheet: new row 5 ---------------------------------
$wsLanguage.Range("A5").Value = "English"
$wsLanguage.Range("B5").Value = "Chinese"

# --- Education sheet: new row 3 ----------------------------------------
$wsEducation.Range("A3").Value = "AIT"
$wsEducation.Range("B3").Value = "ME"
$wsEducation.Range("C3").Value = "GIT"
$wsEducation.Range("D3").Value = "BE"

# --- column widths tweaked on the Language sheet -----------------------
$wsLanguage.Columns.Item(1).ColumnWidth = 14.45
$wsLanguage.Columns.Item(2).ColumnWidth = 11.6

# ---------------------------------------------------------------------
# Leftover cursor/selection position on each touched sheet, matching
# where the editor ended up after typing the new rows.
# ---------------------------------------------------------------------
$wsEducation.Range("I14").Select()
$wsCertificate.Range("G18").Select()
$wsSkill.Range("D17").Select()

# Language is the sheet left active/selected when the file was saved,
# so it's activated (and its selection set) last.
$wsLanguage.Range("E15").Select()
$wsLanguage.Activate()
